$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the existing "IP" header (H1) onto the two new
# header cells so they pick up the same style (bold font, border,
# centered/top alignment) used by the other headers.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for columns I (I0) and J (IF)
$values = @{
    2  = @(7, 7)
    3  = @(1, 4)
    4  = @(1, 6)
    5  = @(1, 5)
    6  = @(1, 7)
    7  = @(1, 6)
    8  = @(1, 6)
    9  = @(1, 5)
    10 = @(1, 3)
    11 = @(1, 5)
    12 = @(1, 3)
    13 = @(1, 1)
    14 = @(1, 6)
    15 = @(1, 8)
    16 = @(1, 5)
    17 = @(1, 6)
    18 = @(1, 6)
    19 = @(1, 8)
    20 = @(1, 7)
    21 = @(1, 2)
    22 = @(1, 6)
    23 = @(1, 7)
    24 = @(1, 7)
    25 = @(1, 5)
    26 = @(1, 8)
    27 = @(1, 6)
    28 = @(1, 7)
    29 = @(1, 5)
    30 = @(1, 3)
    31 = @(1, 3)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
